# Update "想去人数" (F column) values on the 展览, 演出 and 全部类型 sheets
# to reflect the latest generated output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 288
$ws1.Range("F3").Value = 201
$ws1.Range("F4").Value = 2364
$ws1.Range("F5").Value = 1760
$ws1.Range("F6").Value = 339
$ws1.Range("F7").Value = 97
$ws1.Range("F8").Value = 812
$ws1.Range("F9").Value = 165

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 24

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 288
$ws4.Range("F3").Value = 201
$ws4.Range("F4").Value = 2364
$ws4.Range("F5").Value = 1760
$ws4.Range("F6").Value = 339
$ws4.Range("F7").Value = 24
$ws4.Range("F8").Value = 97
$ws4.Range("F9").Value = 812
$ws4.Range("F10").Value = 165
